$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on numeric-looking Price cells so COM keeps them as
# literal strings (matching the workbook's inline-string-typed Price column)
# instead of silently coercing to Double and losing exact formatting
# (trailing zeros, sig figs) on save.
$ws.Range("D2").Value = "67.982.36"
$ws.Range("E2").Value = "  -1.01%  "
$ws.Range("D3").Value = "3.817.54"
$ws.Range("E3").Value = "  -2.47%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "599.64"
$ws.Range("E5").Value = "  -0.77%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.43"
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D7").Value = "3.816.30"
$ws.Range("E7").Value = "  -2.44%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  -0.37%  "
$ws.Range("E10").Value = "  -0.90%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.53"
$ws.Range("E11").Value = "  +1.29%  "
$ws.Range("E12").Value = "  +0.43%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000278"
$ws.Range("E13").Value = "  +9.32%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.02"
$ws.Range("E14").Value = "  -0.73%  "
$ws.Range("D15").Value = "4.457.69"
$ws.Range("E15").Value = "  -2.44%  "
$ws.Range("D16").Value = "3.812.19"
$ws.Range("E16").Value = "  -2.72%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.62"
$ws.Range("E17").Value = "  +2.56%  "
$ws.Range("D18").Value = "67.986.36"
$ws.Range("E18").Value = "  -0.98%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.45"
$ws.Range("E19").Value = "  -0.06%  "
$ws.Range("E20").Value = "  -0.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.87"
$ws.Range("E21").Value = "  -0.08%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "472.21"
$ws.Range("E22").Value = "  -0.43%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.740"
$ws.Range("E23").Value = "  -0.43%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000153"
$ws.Range("E24").Value = "  -9.28%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.63"
$ws.Range("E25").Value = "  -0.29%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.32"
$ws.Range("E26").Value = "  +2.66%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.24"
$ws.Range("E27").Value = "  -0.14%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.35"
$ws.Range("E28").Value = "  +2.89%  "
$ws.Range("E29").Value = "  -0.19%  "
$ws.Range("E30").Value = "  -1.78%  "
$ws.Range("D31").Value = "3.963.77"
$ws.Range("E31").Value = "  -2.47%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.76"
$ws.Range("E32").Value = "  -1.17%  "
$ws.Range("E33").Value = "  -1.59%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "30.84"
$ws.Range("E34").Value = "  -2.70%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.39"
$ws.Range("E35").Value = "  -1.33%  "
$ws.Range("D36").Value = "3.782.86"
$ws.Range("E36").Value = "  -2.67%  "
$ws.Range("E37").Value = "  +1.61%  "
$ws.Range("E38").Value = "  +3.89%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.01"
$ws.Range("E39").Value = "  +0.87%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.139"
$ws.Range("E40").Value = "  -1.23%  "
$ws.Range("E41").Value = "  -2.49%  "
$ws.Range("E42").Value = "  -0.20%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.320"
$ws.Range("E43").Value = "  +1.32%  "
$ws.Range("B44").Value = "Cosmos"
$ws.Range("C44").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.83"
$ws.Range("E44").Value = "  +1.99%  "
$ws.Range("B45").Value = "USDe"
$ws.Range("C45").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.98"
$ws.Range("E46").Value = "  -1.89%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "412.97"
$ws.Range("E47").Value = "  -4.27%  "
$ws.Range("B48").Value = "FLOKI"
$ws.Range("C48").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.000288"
$ws.Range("E48").Value = "  -4.77%  "
$ws.Range("B49").Value = "OKB"
$ws.Range("C49").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "46.53"
$ws.Range("E49").Value = "  -1.60%  "
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "142.87"
$ws.Range("E50").Value = "  -0.66%  "
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0361"
$ws.Range("E51").Value = "  -0.03%  "
